$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 27
$ws.Range("A3").Value = 49
$ws.Range("A5").Value = 56
$ws.Range("A6").Value = 35
$ws.Range("A8").Value = 37
$ws.Range("E8").Value = 473
$ws.Range("A14").Value = 46
$ws.Range("A17").Value = 58
$ws.Range("A18").Value = 40
$ws.Range("A19").Value = 30
$ws.Range("E20").Value = 528
$ws.Range("A23").Value = 55
$ws.Range("A24").Value = 9
$ws.Range("A25").Value = 8
$ws.Range("A26").Value = 34
$ws.Range("A27").Value = 57
$ws.Range("A28").Value = 16
$ws.Range("A33").Value = 31
$ws.Range("A34").Value = 51
$ws.Range("A35").Value = 7
$ws.Range("A37").Value = 39
$ws.Range("A38").Value = 47
$ws.Range("A40").Value = 29
$ws.Range("A41").Value = 38
$ws.Range("A43").Value = 33
$ws.Range("A44").Value = 36
$ws.Range("D44").Value = 4.8
$ws.Range("E44").Value = 25
$ws.Range("A47").Value = 54
$ws.Range("A48").Value = 15
$ws.Range("E48").Value = 845
$ws.Range("A55").Value = 28
$ws.Range("A59").Value = 6
